# Fixing errors in example upload files.
$wb = $excel.ActiveWorkbook

# --- "Service Contacts" sheet ---------------------------------------------
$wsSC = $wb.Sheets.Item("Service Contacts")

# New narrower column A so the organisation_path values are readable.
$wsSC.Columns.Item(1).ColumnWidth = 13.7

# Move the active selection from the whole column D to a single cell (D3).
$wsSC.Range("D3").Select()

# --- "Practitioners" sheet --------------------------------------------------
$wsP = $wb.Sheets.Item("Practitioners")
$wsP.Activate()

# Widen a few columns for readability.
$wsP.Columns.Item(1).ColumnWidth = 13.9
$wsP.Columns.Item(3).ColumnWidth = 12.2
$wsP.Columns.Item(6).ColumnWidth = 12.05

# Add the missing practitioner record (row 6) that was dropped from the
# example upload file.
$wsP.Range("A6").Value = "PHN999:NFP02"
$wsP.Range("B6").Value = "P01"
$wsP.Range("C6").Value = 8
$wsP.Range("D6").Value = 1
$wsP.Range("E6").Value = 1973
$wsP.Range("F6").Value = 2
$wsP.Range("G6").Value = 1
$wsP.Range("H6").Value = 1
$wsP.Range("I6").Value = "tag1"

# Reset the active selection on this (now active) sheet to column G.
$wsP.Range("G1:G1048576").Select()
